$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last status check" timestamp in F1 (14:00 -> 14:15)
$ws.Range("F1").Value = "Last status check on: 25.02.2022 14:15"

# 2. D4: change from text "+1.0" to numeric value 1, default style
$ws.Range("D4").Value = 1
$ws.Range("D4").Style = "Normal"

# 3. E4: change from text "2022-02-25 14:00:11" to numeric date serial value
#    44617.58346064815, using the same date number format as the other rows
#    (e.g. E3), which is "YYYY-MM-DD HH:MM:SS".
$ws.Range("E4").NumberFormat = $ws.Range("E3").NumberFormat
$ws.Range("E4").Value = 44617.58346064815
